$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '"A proprietary extract from the echinacea plant (Echinacea purpurea) enhances systemic immune response during a common cold"'
$ws.Range("E2").Value = '[Vinti%Goel%xref no email%1, Ray%Lovlin%xref no email%1, Chuck%Chang%xref no email%1, Jan V.%Slama%xref no email%1, Richard%Barton%xref no email%1, Roland%Gahler%xref no email%1, R.%Bauer%xref no email%1, L.%Goonewardene%xref no email%1, Tapan K.%Basu%xref no email%1]'
$ws.Range("F2").Value = '10.1002/ptr.1733'
$ws.Range("G2").Value = 'CROSSREF'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '2023-06-01'
$ws.Range("H2").ClearFormats()
$ws.Range("E3").Value = '[Steven H.%Yale%xref no email%1,            Kejian%Liu%xref no email%1]'
$ws.Range("I3").Value = '_CROSSREF'
$ws.Range("E4").Value = '[Barret%BP%coreGivesNoEmail%1,                  Bartels%CL%coreGivesNoEmail%1,                  Bauer%R%coreGivesNoEmail%1,                  Brevoort%P%coreGivesNoEmail%1,                  Brinkeborn%RM%coreGivesNoEmail%1,                  Carr%RJ%coreGivesNoEmail%1,                  Ernst%E%coreGivesNoEmail%1,                  Giles%JT%coreGivesNoEmail%1,                  Hoheisel%D%coreGivesNoEmail%1,                  Management%of Influenza in the Southern Hemisphere Trialists Study Group%coreGivesNoEmail%1,                  Muller-Jakic%B%coreGivesNoEmail%1,                  Stimpel%M%coreGivesNoEmail%1,                  Turner%RB%coreGivesNoEmail%1,                  Winther%B%coreGivesNoEmail%1]'
$ws.Range("I4").Value = '_CROSSREF'
$ws.Range("C5").Value = '"Efficacy of Echinacea purpurea in Patients with a Common Cold"'
$ws.Range("E5").Value = '[Brigitte%Schulten%xref no email%1, Michael%Bulitta%xref no email%1, Brigitta%Ballering-Br\u00fchl%xref no email%1, Ulrike%K\u00f6ster%xref no email%1, Michael%Sch\u00e4fer%xref no email%1]'
$ws.Range("F5").Value = '10.1055/s-0031-1300080'
$ws.Range("G5").Value = 'CROSSREF'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '2023-04-25'
$ws.Range("H5").ClearFormats()
$ws.Range("I6").Value = '_CROSSREF'
$ws.Range("E7").Value = '[G. Frank%Lindenmuth%xref no email%1,            Elise B.%Lindenmuth%xref no email%1]'
$ws.Range("I7").Value = '_CROSSREF'
$ws.Range("E8").Value = '[ Wolfram%Grimm%null%1,                    Hans-Helge%Müller%null%1,                  Wolfram%Grimm%null%1,                  Hans-Helge%Müller%null%1]'
$ws.Range("I8").Value = '_CROSSREF'
$ws.Range("I9").Value = '_CROSSREF'
$ws.Range("E10").Value = '[H.%Hall%xref no email%1,            M.%Fahlman%xref no email%1,            H.%Engels%xref no email%1]'
$ws.Range("I10").Value = '_CROSSREF'
$ws.Range("E11").Value = '[Joelle%O’Neil%NULL%1,                   Susan%Hughes%susan.hughes@fresno.ucsf.edu%1,                   Andrea%Lourie%NULL%1,                   John%Zweifler%NULL%1]'
$ws.Range("I11").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E12").Value = '[M.%Jawad%NULL%1,                   R.%Schoop%NULL%1,                   A.%Suter%NULL%1,                   P.%Klein%NULL%1,                   R.%Eccles%NULL%1]'
$ws.Range("I12").Value = '_PMC_CROSSREF'
$ws.Range("E13").Value = '[E.%Tiralongo%NULL%1,                   R. A.%Lea%NULL%1,                   S. S.%Wee%NULL%1,                   M. M.%Hanna%NULL%1,                   L. R.%Griffiths%NULL%1]'
$ws.Range("I13").Value = '_PMC_CROSSREF'
$ws.Range("E14").Value = '[Ronald B.%Turner%xref no email%1,            Rudolf%Bauer%xref no email%1,            Karin%Woelkart%xref no email%1,            Thomas C.%Hulsey%xref no email%1,            J. David%Gangemi%xref no email%1]'
$ws.Range("I14").Value = '_CROSSREF'
$ws.Range("E15").Value = '[Steven J.%Sperber%ssperber@humed.com%1,                   Leena P.%Shah%NULL%1,                   Richard D.%Gilbert%NULL%1,                   Thomas W.%Ritchey%NULL%1,                   Arnold S.%Monto%NULL%1]'
$ws.Range("I15").Value = '_PMC_CROSSREF'
$ws.Range("C16").Value = '"Echinacea purpurea along with zinc, selenium and vitamin C to alleviate exacerbations of chronic obstructive pulmonary disease: results from a randomized controlled trial"'
$ws.Range("E16").Value = '[F.%Isbaniah%xref no email%1, W. H.%Wiyono%xref no email%1, F.%Yunus%xref no email%1, A.%Setiawati%xref no email%1, U.%Totzke%xref no email%1, M. A.%Verbruggen%xref no email%1]'
$ws.Range("F16").Value = '10.1111/j.1365-2710.2010.01212.x'
$ws.Range("G16").Value = 'CROSSREF'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '2023-06-01'
$ws.Range("H16").ClearFormats()
$ws.Range("I17").Value = '_CROSSREF'
$ws.Range("E18").Value = '[ M.%Dorn%null%1,                    E.%Knick%null%1,                    G.%Lewith%null%1,                  M.%Dorn%null%1,                  E.%Knick%null%1,                  G.%Lewith%null%1]'
$ws.Range("I18").Value = '_CROSSREF'
